$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 338.6742503333333
$ws.Range("H2").Value = 1016.022751
$ws.Range("I2").Value = 0.5849329800180821
$ws.Range("J2").Value = 0.584932980018082
$ws.Range("M2").Value = 209.26237
$ws.Range("N2").Value = 627.78711
$ws.Range("O2").Value = 0.8127157202241573
$ws.Range("P2").Value = 0.8127157202241573
$ws.Range("Q2").Value = 70871.77628272661
$ws.Range("R2").Value = 637845.9865445396
$ws.Range("S2").Value = 0.4753842281382582
$ws.Range("T2").Value = 0.4753842281382581
$ws.Range("G3").Value = 338.6742503333333
$ws.Range("H3").Value = 1016.022751
$ws.Range("I3").Value = 0.5849329800180821
$ws.Range("J3").Value = 0.584932980018082
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.003824977881910862
$ws.Range("P3").Value = 0.003824977881910862
$ws.Range("Q3").Value = 333.5520280798752
$ws.Range("R3").Value = 3001.968252718877
$ws.Range("S3").Value = 0.002237355710969373
$ws.Range("T3").Value = 0.002237355710969372
$ws.Range("G4").Value = 338.6742503333333
$ws.Range("H4").Value = 1016.022751
$ws.Range("I4").Value = 0.5849329800180821
$ws.Range("J4").Value = 0.584932980018082
$ws.Range("M4").Value = 1.763846666666667
$ws.Range("N4").Value = 5.291539999999999
$ws.Range("O4").Value = 0.006850280411451801
$ws.Range("P4").Value = 0.006850280411451801
$ws.Range("Q4").Value = 597.3694475362821
$ws.Range("R4").Value = 5376.32502782654
$ws.Range("S4").Value = 0.004006954935029995
$ws.Range("T4").Value = 0.004006954935029995
$ws.Range("G5").Value = 338.6742503333333
$ws.Range("H5").Value = 1016.022751
$ws.Range("I5").Value = 0.5849329800180821
$ws.Range("J5").Value = 0.584932980018082
$ws.Range("M5").Value = 45.474231
$ws.Range("N5").Value = 136.422693
$ws.Range("O5").Value = 0.1766090214824801
$ws.Range("P5").Value = 0.1766090214824801
$ws.Range("Q5").Value = 15400.95109340983
$ws.Range("R5").Value = 138608.5598406884
$ws.Range("S5").Value = 0.1033044412338246
$ws.Range("T5").Value = 0.1033044412338245
$ws.Range("I6").Value = 0.279688040971731
$ws.Range("J6").Value = 0.2796880409717309
$ws.Range("M6").Value = 209.26237
$ws.Range("N6").Value = 627.78711
$ws.Range("O6").Value = 0.8127157202241573
$ws.Range("P6").Value = 0.8127157202241573
$ws.Range("Q6").Value = 33887.6229343229
$ws.Range("R6").Value = 304988.6064089061
$ws.Range("S6").Value = 0.227306867656424
$ws.Range("T6").Value = 0.2273068676564239
$ws.Range("I7").Value = 0.279688040971731
$ws.Range("J7").Value = 0.2796880409717309
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.003824977881910862
$ws.Range("P7").Value = 0.003824977881910862
$ws.Range("S7").Value = 0.00106980057055185
$ws.Range("T7").Value = 0.00106980057055185
$ws.Range("I8").Value = 0.279688040971731
$ws.Range("J8").Value = 0.2796880409717309
$ws.Range("M8").Value = 1.763846666666667
$ws.Range("N8").Value = 5.291539999999999
$ws.Range("O8").Value = 0.006850280411451801
$ws.Range("P8").Value = 0.006850280411451801
$ws.Range("Q8").Value = 285.6345875943311
$ws.Range("R8").Value = 2570.71128834898
$ws.Range("S8").Value = 0.001915941508385978
$ws.Range("T8").Value = 0.001915941508385977
$ws.Range("I9").Value = 0.279688040971731
$ws.Range("J9").Value = 0.2796880409717309
$ws.Range("M9").Value = 45.474231
$ws.Range("N9").Value = 136.422693
$ws.Range("O9").Value = 0.1766090214824801
$ws.Range("P9").Value = 0.1766090214824801
$ws.Range("Q9").Value = 7364.026286026949
$ws.Range("R9").Value = 66276.23657424253
$ws.Range("S9").Value = 0.04939543123636921
$ws.Range("T9").Value = 0.0493954312363692
$ws.Range("G10").Value = 77.79536166666666
$ws.Range("H10").Value = 233.386085
$ws.Range("I10").Value = 0.1343623634996766
$ws.Range("J10").Value = 0.1343623634996766
$ws.Range("M10").Value = 209.26237
$ws.Range("N10").Value = 627.78711
$ws.Range("O10").Value = 0.8127157202241573
$ws.Range("P10").Value = 0.8127157202241573
$ws.Range("Q10").Value = 16279.64175737382
$ws.Range("R10").Value = 146516.7758163644
$ws.Range("S10").Value = 0.1091984050226597
$ws.Range("T10").Value = 0.1091984050226597
$ws.Range("G11").Value = 77.79536166666666
$ws.Range("H11").Value = 233.386085
$ws.Range("I11").Value = 0.1343623634996766
$ws.Range("J11").Value = 0.1343623634996766
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.003824977881910862
$ws.Range("P11").Value = 0.003824977881910862
$ws.Range("Q11").Value = 76.61875868503279
$ws.Range("R11").Value = 689.5688281652951
$ws.Range("S11").Value = 0.0005139330685475304
$ws.Range("T11").Value = 0.0005139330685475303
$ws.Range("G12").Value = 77.79536166666666
$ws.Range("H12").Value = 233.386085
$ws.Range("I12").Value = 0.1343623634996766
$ws.Range("J12").Value = 0.1343623634996766
$ws.Range("M12").Value = 1.763846666666667
$ws.Range("N12").Value = 5.291539999999999
$ws.Range("O12").Value = 0.006850280411451801
$ws.Range("P12").Value = 0.006850280411451801
$ws.Range("Q12").Value = 137.2190893578778
$ws.Range("R12").Value = 1234.9718042209
$ws.Range("S12").Value = 0.0009204198667182011
$ws.Range("T12").Value = 0.000920419866718201
$ws.Range("G13").Value = 77.79536166666666
$ws.Range("H13").Value = 233.386085
$ws.Range("I13").Value = 0.1343623634996766
$ws.Range("J13").Value = 0.1343623634996766
$ws.Range("M13").Value = 45.474231
$ws.Range("N13").Value = 136.422693
$ws.Range("O13").Value = 0.1766090214824801
$ws.Range("P13").Value = 0.1766090214824801
$ws.Range("Q13").Value = 3537.684247158545
$ws.Range("R13").Value = 31839.1582244269
$ws.Range("S13").Value = 0.02372960554175119
$ws.Range("T13").Value = 0.02372960554175118
$ws.Range("G14").Value = 0.5886170000000001
$ws.Range("H14").Value = 1.765851
$ws.Range("I14").Value = 0.001016615510510267
$ws.Range("J14").Value = 0.001016615510510266
$ws.Range("M14").Value = 209.26237
$ws.Range("N14").Value = 627.78711
$ws.Range("O14").Value = 0.8127157202241573
$ws.Range("P14").Value = 0.8127157202241573
$ws.Range("Q14").Value = 123.17538844229
$ws.Range("R14").Value = 1108.57849598061
$ws.Range("S14").Value = 0.0008262194068154007
$ws.Range("T14").Value = 0.0008262194068154005
$ws.Range("G15").Value = 0.5886170000000001
$ws.Range("H15").Value = 1.765851
$ws.Range("I15").Value = 0.001016615510510267
$ws.Range("J15").Value = 0.001016615510510266
$ws.Range("M15").Value = 0.9848756666666668
$ws.Range("N15").Value = 2.954627
$ws.Range("O15").Value = 0.003824977881910862
$ws.Range("P15").Value = 0.003824977881910862
$ws.Range("Q15").Value = 0.5797145602863335
$ws.Range("R15").Value = 5.217431042577001
$ws.Range("S15").Value = 0.000003888531842109289
$ws.Range("T15").Value = 0.000003888531842109289
$ws.Range("G16").Value = 0.5886170000000001
$ws.Range("H16").Value = 1.765851
$ws.Range("I16").Value = 0.001016615510510267
$ws.Range("J16").Value = 0.001016615510510266
$ws.Range("M16").Value = 1.763846666666667
$ws.Range("N16").Value = 5.291539999999999
$ws.Range("O16").Value = 0.006850280411451801
$ws.Range("P16").Value = 0.006850280411451801
$ws.Range("Q16").Value = 1.038230133393333
$ws.Range("R16").Value = 9.34407120054
$ws.Range("S16").Value = 0.000006964101317626552
$ws.Range("T16").Value = 0.000006964101317626551
$ws.Range("G17").Value = 0.5886170000000001
$ws.Range("H17").Value = 1.765851
$ws.Range("I17").Value = 0.001016615510510267
$ws.Range("J17").Value = 0.001016615510510266
$ws.Range("M17").Value = 45.474231
$ws.Range("N17").Value = 136.422693
$ws.Range("O17").Value = 0.1766090214824801
$ws.Range("P17").Value = 0.1766090214824801
$ws.Range("Q17").Value = 26.766905428527
$ws.Range("R17").Value = 240.902148856743
$ws.Range("S17").Value = 0.0001795434705351301
$ws.Range("T17").Value = 0.0001795434705351301
